$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 897.88
$ws.Range("I15").Value = 897.88
$ws.Range("K15").Value = 2693.64
$ws.Range("M15").Value = -2524.64
$ws.Range("H59").Value = 1250
$ws.Range("I59").Value = 1000
$ws.Range("K59").Value = 3000
$ws.Range("M59").Value = -2443
$ws.Range("H63").Value = 49700
$ws.Range("J63").Value = 49700
$ws.Range("L63").Value = 49700
$ws.Range("N63").Value = -50948
$ws.Range("H66").Value = 49700
$ws.Range("J66").Value = 49700
$ws.Range("L66").Value = 149100
$ws.Range("N66").Value = -155340
$ws.Range("H129").Value = 1698.9756
$ws.Range("I129").Value = 584.53845
$ws.Range("J129").Value = 2216.3928
$ws.Range("K129").Value = 1753.61535
$ws.Range("L129").Value = 6649.178400000001
$ws.Range("M129").Value = 3246.38465
$ws.Range("N129").Value = -16649.1784
$ws.Range("H135").Value = 1195.0605
$ws.Range("I135").Value = 935.4091
$ws.Range("J135").Value = 1714.3636
$ws.Range("K135").Value = 8418.6819
$ws.Range("L135").Value = 15429.2724
$ws.Range("M135").Value = -5883.6819
$ws.Range("N135").Value = -20499.2724
$ws.Range("H137").Value = 1333.9517
$ws.Range("I137").Value = 1707.2069
$ws.Range("J137").Value = 1005.9394
$ws.Range("K137").Value = 5121.620699999999
$ws.Range("L137").Value = 3017.8182
$ws.Range("M137").Value = -2571.620699999999
$ws.Range("N137").Value = -8117.8182
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16597.678
$ws.Range("I32").Value = 9536.959000000001
$ws.Range("J32").Value = 49253.5
$ws.Range("K32").Value = 9536.959000000001
$ws.Range("L32").Value = 49253.5
$ws.Range("M32").Value = -9249.959000000001
$ws.Range("N32").Value = -49827.5
$ws.Range("H45").Value = 1485
$ws.Range("I45").Value = 1485
$ws.Range("K45").Value = 1485
$ws.Range("M45").Value = -1108
$ws.Range("H63").Value = 3854.9048
$ws.Range("I63").Value = 2283.7812
$ws.Range("J63").Value = 8882.5
$ws.Range("K63").Value = 2283.7812
$ws.Range("L63").Value = 8882.5
$ws.Range("M63").Value = -1597.7812
$ws.Range("N63").Value = -10254.5
$ws.Range("H66").Value = 3854.9048
$ws.Range("I66").Value = 2283.7812
$ws.Range("J66").Value = 8882.5
$ws.Range("K66").Value = 11418.906
$ws.Range("L66").Value = 44412.5
$ws.Range("M66").Value = -7986.905999999999
$ws.Range("N66").Value = -51276.5
$ws.Range("H111").Value = 45720
$ws.Range("J111").Value = 45720
$ws.Range("L111").Value = 45720
$ws.Range("N111").Value = -53900
$ws.Range("H132").Value = 2021673.8
$ws.Range("I132").Value = 6503.7407
$ws.Range("J132").Value = 3576233.5
$ws.Range("K132").Value = 19511.2221
$ws.Range("L132").Value = 10728700.5
$ws.Range("M132").Value = -16981.2221
$ws.Range("N132").Value = -10733760.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2739.8
$ws.Range("I20").Value = 1422
$ws.Range("J20").Value = 4716.5
$ws.Range("K20").Value = 1422
$ws.Range("L20").Value = 4716.5
$ws.Range("M20").Value = -1175
$ws.Range("N20").Value = -5210.5
$ws.Range("H54").Value = 8197.333000000001
$ws.Range("I54").Value = 4500
$ws.Range("J54").Value = 10046
$ws.Range("K54").Value = 4500
$ws.Range("L54").Value = 10046
$ws.Range("M54").Value = -4016
$ws.Range("N54").Value = -11014
$ws.Range("H86").Value = 1238.2858
$ws.Range("I86").Value = 1073.25
$ws.Range("J86").Value = 1458.3334
$ws.Range("K86").Value = 1073.25
$ws.Range("L86").Value = 1458.3334
$ws.Range("M86").Value = 49.75
$ws.Range("N86").Value = -3704.3334
$ws.Range("H89").Value = 1238.2858
$ws.Range("I89").Value = 1073.25
$ws.Range("J89").Value = 1458.3334
$ws.Range("K89").Value = 5366.25
$ws.Range("L89").Value = 7291.666999999999
$ws.Range("M89").Value = 249.75
$ws.Range("N89").Value = -18523.667
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 10062.077
$ws.Range("I62").Value = 2597
$ws.Range("J62").Value = 14727.75
$ws.Range("K62").Value = 2597
$ws.Range("L62").Value = 14727.75
$ws.Range("M62").Value = -1973
$ws.Range("N62").Value = -15975.75
$ws.Range("H65").Value = 10062.077
$ws.Range("I65").Value = 2597
$ws.Range("J65").Value = 14727.75
$ws.Range("K65").Value = 12985
$ws.Range("L65").Value = 73638.75
$ws.Range("M65").Value = -9865
$ws.Range("N65").Value = -79878.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 30.166666
$ws.Range("I2").Value = 43.857143
$ws.Range("J2").Value = 24.529411
$ws.Range("K2").Value = 43.857143
$ws.Range("L2").Value = 24.529411
$ws.Range("M2").Value = 69.14285699999999
$ws.Range("N2").Value = -250.529411
$ws.Range("H9").Value = 912.7
$ws.Range("I9").Value = 591.3333
$ws.Range("J9").Value = 1394.75
$ws.Range("K9").Value = 591.3333
$ws.Range("L9").Value = 1394.75
$ws.Range("M9").Value = -421.3333
$ws.Range("N9").Value = -1734.75
$ws.Range("H43").Value = 8672.333000000001
$ws.Range("I43").Value = 2017
$ws.Range("J43").Value = 12000
$ws.Range("K43").Value = 2017
$ws.Range("L43").Value = 12000
$ws.Range("M43").Value = -1866
$ws.Range("N43").Value = -12302
$ws.Range("H46").Value = 6614.2856
$ws.Range("I46").Value = 5123.077
$ws.Range("J46").Value = 26000
$ws.Range("K46").Value = 5123.077
$ws.Range("L46").Value = 26000
$ws.Range("M46").Value = -4967.077
$ws.Range("N46").Value = -26312
$ws.Range("H57").Value = 7271.2856
$ws.Range("I57").Value = 2224.75
$ws.Range("J57").Value = 14000
$ws.Range("K57").Value = 2224.75
$ws.Range("L57").Value = 14000
$ws.Range("M57").Value = -1404.75
$ws.Range("N57").Value = -15640
$ws.Range("H80").Value = 2733.5715
$ws.Range("I80").Value = 2693.6875
$ws.Range("J80").Value = 2861.2
$ws.Range("K80").Value = 2693.6875
$ws.Range("L80").Value = 2861.2
$ws.Range("M80").Value = -1695.6875
$ws.Range("N80").Value = -4857.2
$ws.Range("H83").Value = 2733.5715
$ws.Range("I83").Value = 2693.6875
$ws.Range("J83").Value = 2861.2
$ws.Range("K83").Value = 13468.4375
$ws.Range("L83").Value = 14306
$ws.Range("M83").Value = -8476.4375
$ws.Range("N83").Value = -24290
$ws.Range("H132").Value = 4455.6665
$ws.Range("I132").Value = 3382.2
$ws.Range("J132").Value = 4868.5386
$ws.Range("K132").Value = 10146.6
$ws.Range("L132").Value = 14605.6158
$ws.Range("M132").Value = -7616.599999999999
$ws.Range("N132").Value = -19665.6158
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 26318560
$ws.Range("I40").Value = 47621430
$ws.Range("J40").Value = 3247.8235
$ws.Range("K40").Value = 47621430
$ws.Range("L40").Value = 3247.8235
$ws.Range("M40").Value = -47621294
$ws.Range("N40").Value = -3519.8235
$ws.Range("H46").Value = 476703.47
$ws.Range("I46").Value = 502.27274
$ws.Range("J46").Value = 1000524.8
$ws.Range("K46").Value = 502.27274
$ws.Range("L46").Value = 1000524.8
$ws.Range("M46").Value = -314.27274
$ws.Range("N46").Value = -1000900.8
$ws.Range("H55").Value = 114.64706
$ws.Range("I55").Value = 105
$ws.Range("J55").Value = 123.22222
$ws.Range("K55").Value = 105
$ws.Range("L55").Value = 123.22222
$ws.Range("M55").Value = 68
$ws.Range("N55").Value = -469.22222
$ws.Range("H62").Value = 20000
$ws.Range("J62").Value = 20000
$ws.Range("L62").Value = 20000
$ws.Range("N62").Value = -21248
$ws.Range("H65").Value = 20000
$ws.Range("J65").Value = 20000
$ws.Range("L65").Value = 60000
$ws.Range("N65").Value = -66240
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 27833.334
$ws.Range("I51").Value = 22000
$ws.Range("J51").Value = 28363.637
$ws.Range("K51").Value = 22000
$ws.Range("L51").Value = 28363.637
$ws.Range("M51").Value = -21490
$ws.Range("N51").Value = -29383.637
$ws.Range("H52").Value = 6000
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").Value = $null
$ws.Range("H121").Value = 23572.5
$ws.Range("J121").Value = 23572.5
$ws.Range("L121").Value = 23572.5
$ws.Range("N121").Value = -27066.5
$ws.Range("H122").Value = 1485.7059
$ws.Range("I122").Value = 1228.7
$ws.Range("J122").Value = 1852.8572
$ws.Range("K122").Value = 3686.1
$ws.Range("L122").Value = 5558.571599999999
$ws.Range("M122").Value = -1236.1
$ws.Range("N122").Value = -10458.5716
$ws.Range("H132").Value = 2009.6666
$ws.Range("I132").Value = 1683.1852
$ws.Range("J132").Value = 2989.111
$ws.Range("K132").Value = 5049.5556
$ws.Range("L132").Value = 8967.332999999999
$ws.Range("M132").Value = -2519.5556
$ws.Range("N132").Value = -14027.333

Write-Host "Applied all updates"